$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing summary row ("Média", currently row 12) needs to move down to
# row 28 to make room for the new daily rows (2025-04-15 .. 2025-04-30).
$ws.Range("A28").Value = "Média"
$ws.Range("B28").Value = 34

# New daily data for 2025-04-15 through 2025-04-30.
$dates = @(
    "2025-04-15",
    "2025-04-16",
    "2025-04-17",
    "2025-04-18",
    "2025-04-19",
    "2025-04-20",
    "2025-04-21",
    "2025-04-22",
    "2025-04-23",
    "2025-04-24",
    "2025-04-25",
    "2025-04-26",
    "2025-04-27",
    "2025-04-28",
    "2025-04-29",
    "2025-04-30"
)
$totals = @(44, 40, 38, 48, 30, 30, 37, 43, 31, 35, 39, 37, 37, 41, 29, 25)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 12 + $i

    # Force the date-like text to be stored as plain text (matching the
    # other "Data" column cells) instead of letting Excel auto-convert it
    # into a date serial number. Clearing the format afterwards keeps the
    # cell on the default (unstyled) format, same as its neighbours.
    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $dates[$i]
    $ws.Range("A$row").ClearFormats()

    $ws.Range("B$row").Value = $totals[$i]
    $ws.Range("C$row").Value = "-"
}
